$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move "jabatan" column to sit right after "pendidikan_terakhir" (before "lama_bekerja").
# Insert a blank column at E, copy jabatan's values into it (jabatan is pushed to I by
# the insert), then delete the now-empty original jabatan column.
$ws.Columns("E:E").Insert()
$ws.Range("E1:E2").Value = $ws.Range("I1:I2").Value()
$ws.Columns("I:I").Delete()

# Insert a new column for "umur" right after "nama"
$ws.Columns("C:C").Insert()
$ws.Range("C1").Value = "umur"
$ws.Range("C2").Value = 20

# Insert a new column for "nilai_produktivitas" right before "hasil_penilaian_kinerja_sebelumnya"
$ws.Columns("I:I").Insert()
$ws.Range("I1").Value = "nilai_produktivitas"
$ws.Range("I2").Value = 70

# Match column width for the new last column (J) (bestfit-style width in character units)
$ws.Columns("J:J").ColumnWidth = 29.83

# Update selection to match target view state
$ws.Range("J10").Select()
